# Implement excel bridge WIP
# Fill in D1:E1 with the same "de" header value as C1, and D2:E2 with new
# localized values ("Bla" / "Blupp"), extending the used range from A1:C2
# to A1:E2, and moving the selection to E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = $ws.Range("C1").Value()
$ws.Range("E1").Value = $ws.Range("C1").Value()

$ws.Range("D2").Value = "Bla"
$ws.Range("E2").Value = "Blupp"

$ws.Range("E2").Select()
